# Update countries & provincias Spain
# Refresh the COVID-19 "Pais" data sheet:
#  - update the "last updated" timestamp
#  - refresh case numbers for a set of countries
#  - re-rank a handful of countries (Sudan, Sudan del Sur, Nueva Caledonia,
#    Dominica) whose case counts moved them up the sorted (by "Casos
#    totales", descending) list, displacing the countries that used to sit
#    in those rows down by one position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / timestamp -----------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 7 de Mayo de 2020 a las 15:04"

# --- Data rows ----------------------------------------------------------
# Each entry: row, Pais, Casos totales, Nuevos casos, Casos activos,
#             Recuperados, Casos criticos, Muertes hoy, Muertes
$rows = @(
    @(4,   "Estados Unidos",                     1263697, 605,  213109, 975775, 15827, 14, 74813),
    @(20,  "Arabia Saudita",                      33731, 1793,   7798,  25714,   137, 10,   219),
    @(21,  "Suiza",                                30126,   66,  25700,   2618,   121,  3,  1808),
    @(37,  "Rumania",                              14499,  392,   6144,   7474,   234, 17,   881),
    @(57,  "Argentina",                             5208,    0,   1601,   3334,   148,  0,   273),
    @(59,  "Kazajistan",                            4530,  108,   1470,   3030,    31,  0,    30),
    @(73,  "Uzbekistan",                            2266,   33,   1624,    632,     8,  0,    10),
    @(83,  "Republica de Macedonia",                1572,   33,   1079,    404,    21,  1,    89),
    @(94,  "Sudan",                                  930,   78,     92,    786,     0,  3,    52),
    @(95,  "Letonia",                                909,    9,    464,    427,     3,  1,    18),
    @(96,  "Kirguistan",                             895,   24,    637,    246,    13,  0,    12),
    @(97,  "Republica de Chipre",                    883,    0,    296,    572,    15,  0,    15),
    @(98,  "Somalia",                                873,    0,     87,    747,     2,  0,    39),
    @(99,  "Consejo Danes para los Refugiados",      863,   66,    103,    724,     0,  1,    36),
    @(170, "Sudan del Sur",                           74,   16,      0,     74,     0,  0,     0),
    @(171, "Libia",                                   64,    0,     24,     37,     0,  0,     3),
    @(172, "Polinesia Francesa",                      60,    0,     55,      5,     1,  0,     0),
    @(191, "Nueva Caledonia",                         18,    0,     18,      0,     0,  0,     0),
    @(192, "Belice",                                  18,    0,     16,      0,     0,  0,     2),
    @(198, "Dominica",                                16,    0,     14,      2,     0,  0,     0),
    @(199, "Curazao",                                 16,    0,     13,      2,     0,  0,     1)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Range("A$rowNum").Value = $r[1]
    $ws.Range("B$rowNum").Value = $r[2]
    $ws.Range("C$rowNum").Value = $r[3]
    $ws.Range("D$rowNum").Value = $r[4]
    $ws.Range("E$rowNum").Value = $r[5]
    $ws.Range("F$rowNum").Value = $r[6]
    $ws.Range("G$rowNum").Value = $r[7]
    $ws.Range("H$rowNum").Value = $r[8]
}
